# Apply crypto price/volume updates from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.898.29"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "1.903.13"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'0.8030"
$ws.Range("E5").Value = "  +6.05%  "
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.3108"
$ws.Range("E8").Value = "  +2.32%  "
$ws.Range("D9").Value = "'26.31"
$ws.Range("E9").Value = "  +4.08%  "
$ws.Range("D10").Value = "'0.06998"
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("D11").Value = "'0.07989"
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").Value = "1.898.19"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "'0.7374"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").Value = "'5.154"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "'92.14"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "29.899.96"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "'13.91"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "'5.838"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("D19").Value = "'243.86"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").Value = "'0.000007783"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "2.157.28"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "'0.9999"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "'6.875"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Value = "'167.57"
$ws.Range("E25").Value = "  +1.36%  "
$ws.Range("D26").Value = "'9.165"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "'0.1452"
$ws.Range("E27").Value = "  +14.40%  "
$ws.Range("D28").Value = "'18.82"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").Value = "'2.059"
$ws.Range("E29").Value = "  +2.57%  "
$ws.Range("D30").Value = "'1.356"
$ws.Range("E30").Value = "  -1.99%  "
$ws.Range("D31").Value = "'1.511"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "'4.270"
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").Value = "'0.05510"
$ws.Range("E33").Value = "  +4.35%  "
$ws.Range("D34").Value = "'4.047"
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("E35").Value = "  +1.78%  "
$ws.Range("D36").Value = "'0.7295"
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").Value = "'0.01916"
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").Value = "'2.783"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("D40").Value = "'0.4392"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").Value = "'72.00"
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("D42").Value = "'5.955"
$ws.Range("E42").Value = "  -2.97%  "
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "'0.8355"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("D45").Value = "'1.882"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'7.522"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").Value = "'9.669"
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "'978.92"
$ws.Range("E49").Value = "  +8.66%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.062.83"
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").Value = "'36.08"
$ws.Range("E51").Value = "  -0.08%  "
